$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# sheet1 (prep_sheet) content edits
# ---------------------------------------------------------------------------

# New header label mirrored onto the right-hand (D/E) table, matching the
# existing "130pm-3pm" header used lower down in the left-hand (A/B) table.
$ws1.Range("D4").Value = "130pm-3pm"

# "Smoked Salmon Tea Sandwich" block: merge the two prep lines "Cut chives"
# and "Whipped creme fraiche" into a single "Whipped chive creme fraiche"
# line, then shift "Smoked Salmon" up one row and clear the now-unused row.
$ws1.Range("D21").Value = "Whipped chive creme fraiche"
$ws1.Range("D22").Value = "Smoked Salmon"
$ws1.Range("E22").Value = "2x sides"
$ws1.Range("D23:E23").ClearContents()

# Quantity note expanded with an approximate yield callout; row grows taller
# to fit the wrapped text.
$ws1.Range("E17").Value = "6 quarts, *approx 6 airliner brests"
$ws1.Rows.Item(17).RowHeight = 40

# Row 22 no longer needs the taller spacing used while it held the
# (now-removed) wrapped text.
$ws1.Rows.Item(22).RowHeight = 19

# ---------------------------------------------------------------------------
# sheet1 layout / view tweaks
# ---------------------------------------------------------------------------

# New column (F) width used for the wrapped note text.
$ws1.Columns.Item(6).ColumnWidth = 28.666666666666668

# Selection / scroll position left behind by the edit.
$ws1.Range("D28:D31").Select()

# Fit the sheet to one printed page wide/tall at 61% scale.
$ws1.PageSetup.FitToPagesWide = 1
$ws1.PageSetup.FitToPagesTall = 1
$ws1.PageSetup.Zoom = 61

# ---------------------------------------------------------------------------
# sheet2 (order_sheet) layout tweak
# ---------------------------------------------------------------------------
$ws2.PageSetup.Orientation = 1
